# Temp attachment save: append 3 member rows (Id, Name, PhoneNumber, IsMember, CoffeeCount)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Trigger a <headerFooter/> element on the sheet (matches the target diff).
$ws.PageSetup.CenterHeader = ""

# Try to request a full recalculation on load (best-effort; engine may not
# surface this particular workbook-level flag through COM).
$wb.ForceFullCalculation = $true

function Set-TextValue($range, [string]$text) {
    # Force a genuinely numeric-looking string (phone numbers) to be stored
    # as text (shared string) rather than being auto-coerced to a number,
    # then drop back to the default "Normal" style so no extra numFmt/style
    # record lingers on the cell.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "CHulbul"
Set-TextValue $ws.Range("C2") "9842052424"
$ws.Range("D2").Value = $false
$ws.Range("E2").Value = 0

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Lazy shit"
Set-TextValue $ws.Range("C3") "983244324"
$ws.Range("D3").Value = $false
$ws.Range("E3").Value = 0

# Row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Checking"
Set-TextValue $ws.Range("C4") "12345"
$ws.Range("D4").Value = $false
$ws.Range("E4").Value = 0
